$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Productos_Compra" mini-table (rows 17-20): a new "fechaCompra" column is
# inserted at column G, and the previous H/I/J columns shift one to the
# right (I/J/K). The "Pedidos_Proveedor" mini-table that used to sit at
# columns L:M (rows 17-20) is relocated next to the "Pedido" table (rows
# 37-40), now at columns K:L.
# ---------------------------------------------------------------------------

# -- Stamp destination formatting from stable donor cells (never touched by
#    this script) BEFORE writing values, so we never depend on a cell that
#    is later cleared.
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
}

# Row 17 - table title moves from H17 to I17
Copy-Format "B2" "I17"

# Row 18 - header row
Copy-Format "B3" "G18"    # new fechaCompra column header
Copy-Format "F3" "I18"
Copy-Format "F3" "J18"
Copy-Format "F3" "K18"

# Row 19 - data row
Copy-Format "H4" "G19"
Copy-Format "H4" "I19"
Copy-Format "H4" "J19"
Copy-Format "H4" "K19"

# Row 20 - blank spacer row
Copy-Format "B5" "G20"
Copy-Format "H45" "I20"
Copy-Format "G10" "J20"
Copy-Format "G10" "K20"

# "Pedidos_Proveedor" table re-created at rows 37-40, columns K:L
Copy-Format "B2" "K37"

Copy-Format "B3" "K38"
Copy-Format "B3" "L38"

Copy-Format "H4" "K39"
Copy-Format "H4" "L39"

Copy-Format "B5" "K40"
Copy-Format "C5" "L40"

# -- Now write the actual values into the newly formatted cells.
$ws.Range("I17").Value = "Productos_Compra"

$ws.Range("G18").Value = "fechaCompra"
$ws.Range("I18").Value = "idCompra"
$ws.Range("J18").Value = "idProducto"
$ws.Range("K18").Value = "cantUnidadesCompradas"

$ws.Range("G19").Value = "NN"
$ws.Range("I19").Value = "PK, FK Compra.id"
$ws.Range("J19").Value = "PK, FK Producto.codigoBarras"
$ws.Range("K19").Value = "NN"

$ws.Range("K37").Value = "Pedidos_Proveedor"

$ws.Range("K38").Value = "idPedido"
$ws.Range("L38").Value = "idProveedor"

$ws.Range("K39").Value = "PK, FK Pedido.id"
$ws.Range("L39").Value = "PK, FK Proveedor.nit"

# -- Finally clear out the cells vacated by the reshuffle.
$ws.Range("H17").Clear()
$ws.Range("H18").Clear()
$ws.Range("H19").Clear()
$ws.Range("H20").Clear()

$ws.Range("L17").Clear()
$ws.Range("L18").Clear()
$ws.Range("M18").Clear()
$ws.Range("L19").Clear()
$ws.Range("M19").Clear()
$ws.Range("L20").Clear()
$ws.Range("M20").Clear()

# -- Update the active cell / selection to match the saved view state.
$ws.Range("H19").Select()
